$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row of data (row 3) with dates in columns E and F,
# matching the text-formatted style already used by E2:F2.
$ws.Range("E3:F3").NumberFormat = "@"

$ws.Range("E3").Value = "03262022"
$ws.Range("F3").Value = "03292022"

$ws.Range("C3").Select()
